# Update CRMUserPipelineData.xlsx:
#  - add "day"/"time" columns (B, C) with header + 2 data rows to the "event" sheet
#  - apply a time number format to the new time cells
#  - make "event" the active/selected sheet (was "data")
#  - give the "event" sheet a page setup (portrait) like the other sheets

$wb = $excel.ActiveWorkbook

$eventSheet = $wb.Worksheets.Item("event")

# New header cells
$eventSheet.Range("B1").Value = "day"
$eventSheet.Range("C1").Value = "time"

# New data rows
$eventSheet.Range("B2").Value = 16
$eventSheet.Range("C2").Value = 0.54166666666666663
$eventSheet.Range("C2").NumberFormat = "h:mm:ss"

$eventSheet.Range("B3").Value = 17
$eventSheet.Range("C3").Value = 0.70833333333333337
$eventSheet.Range("C3").NumberFormat = "h:mm:ss"

# Page setup, matching the other data sheets
$eventSheet.PageSetup.Orientation = 1

# Switch the active/selected sheet from "data" to "event"
$eventSheet.Range("A1").Select() | Out-Null
$eventSheet.Activate() | Out-Null
